$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "55×12=" "67×38="
Replace-Text "29×32=" "98×66="
Replace-Text "64×49=" "43×63="
Replace-Text "36×17=" "95×48="
Replace-Text "23×52=" "12×74="
Replace-Text "24×47=" "15×52="
Replace-Text "25×30=" "51×32="
Replace-Text "55×26=" "15×51="
Replace-Text "19×24=" "58×52="
Replace-Text "60×95=" "86×22="
Replace-Text "58×34=" "14×25="
Replace-Text "94×83=" "88×47="
Replace-Text "58×93=" "23×48="
Replace-Text "50×83=" "44×50="
Replace-Text "76×34=" "57×79="
Replace-Text "69×45=" "59×68="
Replace-Text "21×98=" "61×47="
Replace-Text "64×55=" "92×19="
Replace-Text "32×20=" "33×22="
Replace-Text "80×44=" "60×23="
Replace-Text "25×82=" "95×95="
Replace-Text "76×24=" "98×26="
Replace-Text "99×37=" "96×74="
Replace-Text "58×91=" "73×59="
Replace-Text "78×64=" "58×51="
